$d = $word.ActiveDocument

# Replace the ellipsis placeholder paragraph with the real sentence about
# version control.
$d.Content.Find.Execute([char]8230, $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "In software development, the version control is the system that is responsible for managing changes to the components or computer programs", `
    2)

# The document ends with a trailing empty paragraph; turn it into the
# "Syed Jafri" line, then append a new paragraph for the GitHub handle.
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$pLast.Range.InsertBefore("Syed Jafri")

$pLast = $d.Paragraphs.Item($lastIndex)
$pLast.Range.InsertParagraphAfter()
$pHandle = $d.Paragraphs.Item($lastIndex + 1)
$pHandle.Range.InsertBefore("github/askarii")
